$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "/account/new/check" right after the "/account/" row
# (columns A:D only, shifting existing rows down)
$ws.Range("A7:D7").Insert(-4121)

# Copy formatting from row 6 (the "/account/" row) into the newly inserted row 7
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the URL for the account registration screen row (row 6)
$ws.Range("B6").Value = "/account/new/"

# Fill in the new confirmation row
$ws.Range("B7").Value = "/account/new/check"
$ws.Range("C7").Value = "登録の確認処理と画面"

# Update the URL for the account registration processing row (now row 8)
$ws.Range("B8").Value = "/account/new/add.php"

# Update selection to match the authored state
$ws.Range("C7").Select()
